$d = $word.ActiveDocument

# 1. "Basic xtUML Modeling (four days)" -> "Basic xtUML Modeling (4-5 days)"
$d.Content.Find.Execute("Basic xtUML Modeling (four days)", $true, $false, $false, $false, $false, $true, 1, $false, "Basic xtUML Modeling (4-5 days)", 2) | Out-Null

# 2. "Tool Training (TBD)" -> "Tool Training (one day)"
$d.Content.Find.Execute("Tool Training (TBD)", $true, $false, $false, $false, $false, $true, 1, $false, "Tool Training (one day)", 2) | Out-Null

# 3. "Completion of Case Study Model (TBD)" -> "Completion of Case Study Model (four days)"
$d.Content.Find.Execute("Completion of Case Study Model (TBD)", $true, $false, $false, $false, $false, $true, 1, $false, "Completion of Case Study Model (four days)", 2) | Out-Null

# 4. Update the "capable of building an xtUML model..." sentence.
$d.Content.Find.Execute("Upon completion of this component the student is capable of building an xtUML model, given only a functional specification as input.  This includes:", $true, $false, $false, $false, $false, $true, 1, $false, "Upon completion of this component the student is capable of building an xtUML model, given a functional specification and requirements-clarification models as input.  This includes:", 2) | Out-Null

# 5. Remove the "Connecting models to externally-produced code" list item (the hyphenated
#    variant in numId=5), located right before "Creating modelled test cases".
$r = $d.Content
$r.Find.Execute("Connecting models to externally-produced code") | Out-Null
$para = $r.Paragraphs(1)
$para.Range.Delete()

# 6. "Given limited time (four days) and the student's finite capacity..." -> "Given the student's finite capacity..."
$d.Content.Find.Execute("Given limited time (four days) and the student", $true, $false, $false, $false, $false, $true, 1, $false, "Given the student", 2) | Out-Null

# 7. Insert a new paragraph describing the workshop days, right after the
#    "...sufficient for building working models." paragraph.
$r2 = $d.Content
$r2.Find.Execute("sufficient for building working models.") | Out-Null
$para2 = $r2.Paragraphs(1)
$idx = $para2.Range.Paragraphs(1).Index
$para2.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($idx + 1)
# A temporary marker is appended so the _GoBack bookmark (below) can be anchored
# with a real, non-boundary position, then the marker is stripped back out.
$newPara.Range.Text = "Days one through three comprise a mixture of lectures and short, focused workshops.  Days four and five are dedicated entirely to workshops during which the students work in small teams of 2-3 each to complete (as much as time permits) the case study model.  An instructor provides consulting and guidance during these workshop days, and the students explain their models to their peers as well as critique the models produced by other teams.ZZGOBACKMARKERZZ"

# Move the _GoBack bookmark to the end of the freshly-inserted paragraph, mirroring the
# position it occupies after this being the most recent edit.
$bmRange = $d.Content
$bmRange.Find.Execute("other teams.") | Out-Null
$zeroRange = $d.Range($bmRange.End, $bmRange.End)
$d.Bookmarks.Add("_GoBack", $zeroRange)

$markerRange = $d.Content
$markerRange.Find.Execute("ZZGOBACKMARKERZZ") | Out-Null
$markerRange.Delete()

Write-Output "done"
